$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously missing values in row 3
$ws.Range("H3").Value = 10042
$ws.Range("I3").Value = 3201

# Move the active selection to J8 (matches the saved selection state in the diff)
$ws.Range("J8").Select()
